$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.797.34'
$ws.Range("E2").Value = '  +1.01%  '

$ws.Range("D3").Value = '3.513.67'
$ws.Range("E3").Value = '  -0.10%  '

$ws.Range("E4").Value = '  +0.05%  '

$c = $ws.Range("D5")
$c.NumberFormat = '@'
$c.Value = '604.75'
$c.Style = 'Normal'
$ws.Range("E5").Value = '  +4.09%  '

$c = $ws.Range("D6")
$c.NumberFormat = '@'
$c.Value = '171.39'
$c.Style = 'Normal'
$ws.Range("E6").Value = '  -2.06%  '

$ws.Range("E7").Value = '  -1.50%  '

$ws.Range("D8").Value = '3.510.29'
$ws.Range("E8").Value = '  -0.04%  '

$ws.Range("E9").Value = '  -0.03%  '

$c = $ws.Range("D10")
$c.NumberFormat = '@'
$c.Value = '0.198'
$c.Style = 'Normal'
$ws.Range("E10").Value = '  +4.30%  '

$ws.Range("E11").Value = '  -1.40%  '

$c = $ws.Range("D12")
$c.NumberFormat = '@'
$c.Value = '0.582'
$c.Style = 'Normal'
$ws.Range("E12").Value = '  -3.04%  '

$c = $ws.Range("D13")
$c.NumberFormat = '@'
$c.Value = '47.23'
$c.Style = 'Normal'
$ws.Range("E13").Value = '  -0.25%  '

$c = $ws.Range("D14")
$c.NumberFormat = '@'
$c.Value = '0.0000279'
$c.Style = 'Normal'
$ws.Range("E14").Value = '  +0.61%  '

$ws.Range("D15").Value = '4.083.75'
$ws.Range("E15").Value = '  +0.04%  '

$c = $ws.Range("D16")
$c.NumberFormat = '@'
$c.Value = '619.66'
$c.Style = 'Normal'
$ws.Range("E16").Value = '  -8.42%  '

$c = $ws.Range("D17")
$c.NumberFormat = '@'
$c.Value = '8.38'
$c.Style = 'Normal'
$ws.Range("E17").Value = '  -4.62%  '

$ws.Range("D18").Value = '3.514.55'
$ws.Range("E18").Value = '  -0.29%  '

$ws.Range("D19").Value = '69.863.79'
$ws.Range("E19").Value = '  +1.15%  '

$ws.Range("E20").Value = '  -2.15%  '

$c = $ws.Range("D21")
$c.NumberFormat = '@'
$c.Value = '17.29'
$c.Style = 'Normal'
$ws.Range("E21").Value = '  -1.66%  '

$c = $ws.Range("D22")
$c.NumberFormat = '@'
$c.Value = '9.97'
$c.Style = 'Normal'
$ws.Range("E22").Value = '  -11.54%  '

$c = $ws.Range("D23")
$c.NumberFormat = '@'
$c.Value = '0.884'
$c.Style = 'Normal'
$ws.Range("E23").Value = '  -2.69%  '

$c = $ws.Range("D24")
$c.NumberFormat = '@'
$c.Value = '15.73'
$c.Style = 'Normal'
$ws.Range("E24").Value = '  -2.89%  '

$c = $ws.Range("D25")
$c.NumberFormat = '@'
$c.Value = '95.93'
$c.Style = 'Normal'
$ws.Range("E25").Value = '  -2.40%  '

$c = $ws.Range("D26")
$c.NumberFormat = '@'
$c.Value = '3.84'
$c.Style = 'Normal'
$ws.Range("E26").Value = '  -0.75%  '

$ws.Range("E27").Value = '  +0.06%  '

$ws.Range("E28").Value = '  -3.04%  '

$c = $ws.Range("D29")
$c.NumberFormat = '@'
$c.Value = '9.22'
$c.Style = 'Normal'
$ws.Range("E29").Value = '  -2.91%  '

$c = $ws.Range("D30")
$c.NumberFormat = '@'
$c.Value = '33.15'
$c.Style = 'Normal'
$ws.Range("E30").Value = '  +0.41%  '

$c = $ws.Range("D31")
$c.NumberFormat = '@'
$c.Value = '8.43'
$c.Style = 'Normal'
$ws.Range("E31").Value = '  -4.17%  '

$c = $ws.Range("D32")
$c.NumberFormat = '@'
$c.Value = '3.07'
$c.Style = 'Normal'
$ws.Range("E32").Value = '  -4.79%  '

$ws.Range("E33").Value = '  -1.69%  '

$c = $ws.Range("D34")
$c.NumberFormat = '@'
$c.Value = '6.97'
$c.Style = 'Normal'
$ws.Range("E34").Value = '  -5.86%  '

$c = $ws.Range("D35")
$c.NumberFormat = '@'
$c.Value = '566.58'
$c.Style = 'Normal'
$ws.Range("E35").Value = '  -2.32%  '

$ws.Range("E36").Value = '  -1.51%  '

$c = $ws.Range("D37")
$c.NumberFormat = '@'
$c.Value = '3.48'
$c.Style = 'Normal'
$ws.Range("E37").Value = '  -3.51%  '

$c = $ws.Range("D38")
$c.NumberFormat = '@'
$c.Value = '57.05'
$c.Style = 'Normal'
$ws.Range("E38").Value = '  -0.62%  '

$ws.Range("E39").Value = '  -3.87%  '

$ws.Range("E40").Value = '  +0.06%  '

$ws.Range("B41").Value = 'Kaspa'
$ws.Range("C41").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$c = $ws.Range("D41")
$c.NumberFormat = '@'
$c.Value = '0.140'
$c.Style = 'Normal'
$ws.Range("E41").Value = '  +2.74%  '

$ws.Range("B42").Value = 'VeChain'
$ws.Range("C42").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$c = $ws.Range("D42")
$c.NumberFormat = '@'
$c.Value = '0.0449'
$c.Style = 'Normal'
$ws.Range("E42").Value = '  +1.67%  '

$ws.Range("E43").Value = '  -3.96%  '

$ws.Range("D44").Value = '3.330.43'
$ws.Range("E44").Value = '  -3.00%  '

$c = $ws.Range("D46")
$c.NumberFormat = '@'
$c.Value = '33.08'
$c.Style = 'Normal'
$ws.Range("E46").Value = '  -1.55%  '

$ws.Range("D47").Value = '0.0₃0703'
$ws.Range("E47").Value = '  -1.12%  '

$ws.Range("E48").Value = '  +0.74%  '

$ws.Range("E49").Value = '  -3.70%  '

$c = $ws.Range("D50")
$c.NumberFormat = '@'
$c.Value = '135.93'
$c.Style = 'Normal'
$ws.Range("E50").Value = '  +2.90%  '

$c = $ws.Range("D51")
$c.NumberFormat = '@'
$c.Value = '5.69'
$c.Style = 'Normal'
$ws.Range("E51").Value = '  +0.61%  '
